$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2023-07-14)
$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 19.60365773276954

# Row 3 (2023-04-03)
$ws.Range("B3").Value = 0.02258322285507441
$ws.Range("C3").Value = 0.004309184025731883
$ws.Range("D3").Value = 2938.103010863317
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 2944.611331347469

# Row 4 (2023-04-01)
$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 5.488907176552729
